$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh crypto price / 1h-volume figures (and swap the last listed coin)
# as scraped on 2024-11-21. Price cells whose new text could be parsed as a
# number are written with a leading apostrophe so Excel keeps them as text,
# matching the original inlineStr cell type.

$ws.Range('D2').Value = '97.011.53'
$ws.Range('E2').Value = '  +4.87%  '
$ws.Range('D3').Value = '3.106.61'
$ws.Range('E3').Value = '  -0.06%  '
$ws.Range('E4').Value = '  +0.23%  '
$ws.Range('D5').Value = "'239.00"
$ws.Range('E5').Value = '  +1.85%  '
$ws.Range('D6').Value = "'609.02"
$ws.Range('E6').Value = '  -0.69%  '
$ws.Range('D7').Value = "'1.11"
$ws.Range('E7').Value = '  +2.39%  '
$ws.Range('D8').Value = "'0.382"
$ws.Range('E8').Value = '  -1.96%  '
$ws.Range('E9').Value = '  +0.21%  '
$ws.Range('D10').Value = '3.095.95'
$ws.Range('E10').Value = '  -0.33%  '
$ws.Range('D11').Value = "'0.784"
$ws.Range('E11').Value = '  +0.22%  '
$ws.Range('E12').Value = '  -0.28%  '
$ws.Range('D13').Value = '96.635.76'
$ws.Range('E13').Value = '  +4.80%  '
$ws.Range('D14').Value = "'0.0000240"
$ws.Range('E14').Value = '  -1.55%  '
$ws.Range('D15').Value = "'33.77"
$ws.Range('E15').Value = '  -0.23%  '
$ws.Range('D16').Value = "'5.35"
$ws.Range('E16').Value = '  -1.45%  '
$ws.Range('D17').Value = '3.699.02'
$ws.Range('E17').Value = '  +0.14%  '
$ws.Range('D18').Value = '3.106.79'
$ws.Range('E18').Value = '  -0.74%  '
$ws.Range('D19').Value = "'3.57"
$ws.Range('E19').Value = '  -6.18%  '
$ws.Range('D20').Value = "'507.73"
$ws.Range('E20').Value = '  +15.60%  '
$ws.Range('D21').Value = "'14.49"
$ws.Range('E21').Value = '  -0.37%  '
$ws.Range('D22').Value = "'5.66"
$ws.Range('E22').Value = '  -2.47%  '
$ws.Range('D23').Value = "'0.0000193"
$ws.Range('E23').Value = '  -5.33%  '
$ws.Range('D24').Value = "'8.80"
$ws.Range('E24').Value = '  -4.71%  '
$ws.Range('D25').Value = "'5.51"
$ws.Range('E25').Value = '  -1.39%  '
$ws.Range('D26').Value = "'86.20"
$ws.Range('E26').Value = '  +1.03%  '
$ws.Range('D27').Value = "'11.61"
$ws.Range('E27').Value = '  +0.77%  '
$ws.Range('D28').Value = '3.279.49'
$ws.Range('E28').Value = '  +0.16%  '
$ws.Range('D29').Value = "'1.00"
$ws.Range('E29').Value = '  +0.19%  '
$ws.Range('D30').Value = "'0.236"
$ws.Range('E30').Value = '  +1.74%  '
$ws.Range('E31').Value = '  -1.22%  '
$ws.Range('D32').Value = "'0.123"
$ws.Range('E32').Value = '  +2.76%  '
$ws.Range('D33').Value = "'8.98"
$ws.Range('E33').Value = '  -2.11%  '
$ws.Range('D34').Value = "'26.49"
$ws.Range('E34').Value = '  +2.38%  '
$ws.Range('D35').Value = "'0.995"
$ws.Range('E35').Value = '  -0.47%  '
$ws.Range('D36').Value = "'0.151"
$ws.Range('E36').Value = '  -3.83%  '
$ws.Range('D37').Value = "'7.32"
$ws.Range('E37').Value = '  -8.91%  '
$ws.Range('D38').Value = "'489.12"
$ws.Range('E38').Value = '  +4.55%  '
$ws.Range('E39').Value = '  -1.27%  '
$ws.Range('D40').Value = "'24.22"
$ws.Range('E40').Value = '  +1.49%  '
$ws.Range('D41').Value = "'0.435"
$ws.Range('E41').Value = '  +0.95%  '
$ws.Range('D42').Value = "'1.24"
$ws.Range('E42').Value = '  -2.91%  '
$ws.Range('D43').Value = "'3.58"
$ws.Range('E43').Value = '  -10.41%  '
$ws.Range('E44').Value = '  +0.02%  '
$ws.Range('D45').Value = "'3.18"
$ws.Range('E45').Value = '  -2.94%  '
$ws.Range('D46').Value = "'162.70"
$ws.Range('E46').Value = '  +1.85%  '
$ws.Range('D47').Value = "'1.90"
$ws.Range('E47').Value = '  +3.26%  '
$ws.Range('D48').Value = "'0.690"
$ws.Range('E48').Value = '  +0.88%  '
$ws.Range('D49').Value = "'44.35"
$ws.Range('E49').Value = '  +1.23%  '
$ws.Range('D50').Value = "'4.35"
$ws.Range('E50').Value = '  +0.08%  '
$ws.Range('B51').Value = 'FirstDigitalUSD'
$ws.Range('C51').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D51').Value = "'1.00"
$ws.Range('E51').Value = '  +0.23%  '
